# Update the cached "today" date shown in the Date placeholders on the
# slide master, the notes master, and every slide layout.
#
#  - notesMaster: "datetimeFigureOut" field   9/9/19   -> 9/18/19
#  - slideMaster + every slideLayout: "datetime1" field 09/09/19 -> 18/09/19

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "09/09/19") {
                $tr.Text = "18/09/19"
            } elseif ($tr.Text -eq "9/9/19") {
                $tr.Text = "9/18/19"
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes

# Every slide layout (CustomLayouts) hanging off the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}
